$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A11").Value = 0
$ws.Range("B11").Value = ""
$ws.Range("C11").Value = ""
$ws.Range("D11").Value = ""
$ws.Range("E11").Value = ""
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = ""
$ws.Range("H11").Value = "..._...@...."
$ws.Range("I11").Value = "**********"
